$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Range("A2").Value = "mohit"
$ws.Range("C2").Value = 2

# Remove rows 3 and 4 entirely (shrinks used range to A1:C2)
$ws.Range("A3:C4").Delete()
